$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update names used in the database testing rows
$ws.Range("A2").Value = "ravi"
$ws.Range("A3").Value = "vipin"

# Update the active view: scroll back to A1 and move the selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("D10").Select()
